# Auto-generated edit script: update crypto price/volume data per commit
# Commit message: Updated cryptos list on Wed Jan  3 07:10:11 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.495.09"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.376.76"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.49"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.83"
$ws.Range("E6").Value = "  -3.48%  "
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.19"
$ws.Range("E10").Value = "  -3.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0921"
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.736.88"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.48"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.360.89"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.427.60"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.41"
$ws.Range("E19").Value = "  +14.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.33"
$ws.Range("E20").Value = "  -3.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000107"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.71"
$ws.Range("E22").Value = "  +5.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.35"
$ws.Range("E23").Value = "  -2.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.30"
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.60"
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.28"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("E30").Value = "  -1.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.68"
$ws.Range("E31").Value = "  -3.64%  "
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "168.57"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("E36").Value = "  -2.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.77"
$ws.Range("E37").Value = "  -3.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.93"
$ws.Range("E38").Value = "  +12.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.02"
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.79"
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("B43").Value = "BitcoinSV"
$ws.Range("C43").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "97.21"
$ws.Range("E43").Value = "  -7.00%  "
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.00"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.229"
$ws.Range("E45").Value = "  -4.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.848.80"
$ws.Range("E46").Value = "  +12.61%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.95"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.07"
$ws.Range("E49").Value = "  +5.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.71"
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.29"
$ws.Range("E51").Value = "  -0.86%  "
